$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_cases")

# Update test case names in column C
$ws.Range("C4").Value = "test_switching_project_list_views"
$ws.Range("C5").Value = "test_signing _out"
$ws.Range("C9").Value = "test_preventing_unsuccessful_login_attempts"

# Apply AutoFilter on the table, filtering column A ("PAGE") to show only "login page"
$ws.Range("A2:H10").AutoFilter(1, @("login page"), 7)

# Update the active selection
$ws.Range("D18").Select()
